$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string rich-text runs) ---
# "Volume 31   Number  9"  ->  "Volume 31   Number  10"
$ws.Range("A8").Value = "Volume 31   Number  10"
# "Report Covering the Week  2/26/2024  Through  3/3/2024"
#   -> "Report Covering the Week  3/4/2024  Through  3/10/2024"
$ws.Range("C9").Value = "Report Covering the Week  3/4/2024  Through  3/10/2024"

# --- Cells that flip between text("0"/blank) and numeric representations ---
# Copy number-format/style from a same-column sibling that already carries the
# target style, then overwrite with the new numeric value, so the style index
# follows the data type exactly like Excel would do on manual entry.
$ws.Range("C15").Copy($ws.Range("C16"))
$ws.Range("C16").Value = 1

$ws.Range("D20").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 7

$ws.Range("F29").Copy($ws.Range("F31"))
$ws.Range("F31").Value = 1

$ws.Range("I29").Copy($ws.Range("I31"))
$ws.Range("I31").Value = 1

# C33 goes from numeric 1 back to the text placeholder "0"
$ws.Range("D27").Copy($ws.Range("C33"))

# --- Plain numeric updates (style unchanged) ---
# Row 15
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 5
$ws.Range("K15").Value = 66.666666666666
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -16.666666666666

# Row 16 (continued)
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -40
$ws.Range("I16").Value = 29
$ws.Range("J16").Value = 26
$ws.Range("K16").Value = 11.538461538461
$ws.Range("L16").Value = -6.451612903225
$ws.Range("M16").Value = -58.571428571428
$ws.Range("N16").Value = -84.895833333333

# Row 17
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -14.285714285714
$ws.Range("F17").Value = 31
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = 24
$ws.Range("I17").Value = 88
$ws.Range("J17").Value = 75
$ws.Range("K17").Value = 17.333333333333
$ws.Range("L17").Value = 8.641975308641
$ws.Range("M17").Value = 66.037735849056
$ws.Range("N17").Value = 35.384615384615

# Row 18
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 9
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -59.259259259259
$ws.Range("I18").Value = 30
$ws.Range("J18").Value = 45
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -6.25
$ws.Range("M18").Value = -55.223880597014
$ws.Range("N18").Value = -87.654320987654

# Row 19
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 37
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = -19.565217391304
$ws.Range("I19").Value = 119
$ws.Range("J19").Value = 114
$ws.Range("K19").Value = 4.385964912280
$ws.Range("L19").Value = 22.680412371134
$ws.Range("M19").Value = 54.545454545454
$ws.Range("N19").Value = 23.958333333333

# Row 20 (continued)
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 133.333333333333
$ws.Range("G20").Value = 27
$ws.Range("H20").Value = -29.629629629629
$ws.Range("I20").Value = 66
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = 26.923076923076
$ws.Range("L20").Value = 73.684210526315
$ws.Range("M20").Value = -5.714285714285
$ws.Range("N20").Value = -89.320388349514

# Row 21 (TOTAL)
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 137
$ws.Range("H21").Value = -21.897810218978
$ws.Range("I21").Value = 338
$ws.Range("J21").Value = 315
$ws.Range("K21").Value = 7.301587301587
$ws.Range("L21").Value = 19.434628975265
$ws.Range("M21").Value = -2.312138728323
$ws.Range("N21").Value = -72.363041700735

# Row 24
$ws.Range("C24").Value = 16
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -60.975609756097
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -27.049180327868
$ws.Range("I24").Value = 227
$ws.Range("J24").Value = 258
$ws.Range("K24").Value = -12.015503875969
$ws.Range("L24").Value = -4.219409282700
$ws.Range("M24").Value = 60.992907801418

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 300
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -11.764705882352
$ws.Range("I25").Value = 45
$ws.Range("J25").Value = 44
$ws.Range("K25").Value = 2.272727272727
$ws.Range("L25").Value = -4.255319148936

# Row 26
$ws.Range("C26").Value = 25
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = 92.307692307692
$ws.Range("F26").Value = 76
$ws.Range("G26").Value = 46
$ws.Range("H26").Value = 65.217391304347
$ws.Range("I26").Value = 138
$ws.Range("J26").Value = 102
$ws.Range("K26").Value = 35.294117647058
$ws.Range("L26").Value = 38
$ws.Range("M26").Value = 7.8125

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 8
$ws.Range("K27").Value = 14.285714285714
$ws.Range("L27").Value = -11.111111111111

# Row 28
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 33.333333333333
$ws.Range("I28").Value = 10
$ws.Range("K28").Value = 42.857142857142
$ws.Range("L28").Value = 66.666666666666

# Row 29
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -77.777777777777

# Row 30
$ws.Range("M30").Value = -50
$ws.Range("N30").Value = -77.777777777777

# Row 31 (continued)
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 1
